$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208 - this shifts the existing row 208 (and everything
# below it) down by one, turning the old row 208 into row 209, etc.
$ws.Rows(208).Insert()

# The new row 208 keeps the same Mercado/Region/Categoria/Calidad/Unidad/
# Origen/Kg-o-Unidades/Clasificacion values as the row that used to be here
# (now shifted to row 209), and only gets fresh Fecha/Volumen/Precio values.
$ws.Range("A208").Value = $ws.Range("A209").Value2
$ws.Range("B208").Value = $ws.Range("B209").Value2
$ws.Range("C208").Value = $ws.Range("C209").Value2
$ws.Range("D208").Value = 44609
$ws.Range("E208").Value = $ws.Range("E209").Value2
$ws.Range("F208").Value = $ws.Range("F209").Value2
$ws.Range("G208").Value = $ws.Range("G209").Value2
$ws.Range("H208").Value = $ws.Range("H209").Value2
$ws.Range("I208").Value = $ws.Range("I209").Value2
$ws.Range("J208").Value = 80
$ws.Range("K208").Value = 5000
$ws.Range("L208").Value = 5000
$ws.Range("M208").Value = 5000
$ws.Range("N208").Value = $ws.Range("N209").Value2
$ws.Range("O208").Value = $ws.Range("O209").Value2
$ws.Range("P208").Value = 2500
$ws.Range("Q208").Value = $ws.Range("Q209").Value2
$ws.Range("R208").Value = $ws.Range("R209").Value2
